$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sessions")

# Rename header cells to avoid the string "file_" in file related settings.
$ws.Range("T1").Value = "fname_1"
$ws.Range("U1").Value = "fposition_1clip_out_1"
$ws.Range("X1").Value = "fname_2fposition_2"
